# fix(publipostage): Try to solve Excel emoji problem
#
# Replace the emoji values used in the "statut" column (A) with
# plain-text / alternate-emoji equivalents:
#   📕 -> -3
#   📘 -> ⚠️
#   📗 -> ✅
#   📙 -> +3

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$map = @{
    "📕" = "-3"
    "📘" = "⚠️"
    "📗" = "✅"
    "📙" = "+3"
}

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row  # xlUp = -4162

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $val = $cell.Value()
    if ($map.ContainsKey($val)) {
        $newVal = $map[$val]
        # "-3"/"+3" look like numbers to Excel's auto-detection, so
        # force text storage first (otherwise they'd be saved as
        # numeric values instead of shared-string text). Restore the
        # default "Normal" style afterwards so no visible formatting
        # change is introduced.
        if ($newVal -eq "-3" -or $newVal -eq "+3") {
            $cell.NumberFormat = "@"
            $cell.Value = $newVal
            $cell.Style = "Normal"
        } else {
            $cell.Value = $newVal
        }
    }
}
